$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column L ("Commune"), which shifts column M ("Indicator Type") into L.
$ws.Range("L:L").Delete()
